$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.440.33'
$ws.Range("E2").Value = '  -1.53%  '
$ws.Range("D3").Value = '1.688.33'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '315.54'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = '0.3903'
$ws.Range("E7").Value = '  -0.80%  '
$ws.Range("D8").Value = '0.4037'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("E9").Value = '  -1.84%  '
$ws.Range("D10").Value = '0.9994'
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("D11").Value = '52.46'
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").Value = '0.08777'
$ws.Range("E12").Value = '  -1.36%  '
$ws.Range("D13").Value = '26.58'
$ws.Range("E13").Value = '  +12.95%  '
$ws.Range("D14").Value = '7.477'
$ws.Range("E14").Value = '  +2.23%  '
$ws.Range("D15").Value = '8.184'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").Value = '0.00001348'
$ws.Range("E16").Value = '  +1.17%  '
$ws.Range("D17").Value = '1.691.23'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '98.08'
$ws.Range("E18").Value = '  -2.13%  '
$ws.Range("D19").Value = '0.07246'
$ws.Range("E19").Value = '  +2.81%  '
$ws.Range("D20").Value = '20.36'
$ws.Range("E20").Value = '  +3.29%  '
$ws.Range("D21").Value = '7.296'
$ws.Range("D22").Value = '0.9997'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = '14.25'
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("D24").Value = '24.432.74'
$ws.Range("E24").Value = '  -1.54%  '
$ws.Range("D25").Value = '3.041'
$ws.Range("E25").Value = '  -4.88%  '
$ws.Range("D26").Value = '2.335'
$ws.Range("E26").Value = '  -1.14%  '
$ws.Range("D27").Value = '22.64'
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("D28").Value = '167.46'
$ws.Range("E28").Value = '  +3.25%  '
$ws.Range("D29").Value = '8.599'
$ws.Range("E29").Value = '  +2.54%  '
$ws.Range("D30").Value = '5.346'
$ws.Range("E30").Value = '  +3.45%  '
$ws.Range("D31").Value = '138.64'
$ws.Range("E31").Value = '  +1.46%  '
$ws.Range("D32").Value = '1.875.05'
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("D33").Value = '0.08780'
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").Value = '7.300'
$ws.Range("E34").Value = '  -2.26%  '
$ws.Range("D35").Value = '2.109'
$ws.Range("E35").Value = '  +6.68%  '
$ws.Range("D36").Value = '1.051'
$ws.Range("E36").Value = '  -3.10%  '
$ws.Range("D37").Value = '0.03029'
$ws.Range("E37").Value = '  +9.80%  '
$ws.Range("D38").Value = '0.2787'
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("E39").Value = '  -3.34%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.09133'
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '14.19'
$ws.Range("E41").Value = '  -1.62%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.8030'
$ws.Range("E42").Value = '  +4.08%  '
$ws.Range("D43").Value = '1.469'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").Value = '17.72'
$ws.Range("E44").Value = '  +11.30%  '
$ws.Range("D45").Value = '2.654'
$ws.Range("E45").Value = '  +2.95%  '
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("D47").Value = '4.263'
$ws.Range("E47").Value = '  +1.34%  '
$ws.Range("D48").Value = '1.427'
$ws.Range("E48").Value = '  +8.03%  '
$ws.Range("D49").Value = '0.9998'
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").Value = '139.39'
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").Value = '0.08071'
$ws.Range("E51").Value = '  +1.01%  '
